# Updated symbol list on Fri Dec 23 20:48:41 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) snapshot for most coins, flips the
# "Worst in 24h" badge from the KickToken row to the One row, and re-sorts
# three rows (BKEXToken / CEJI / KickToken) into their new rank order.
#
# All of column D's values are stored as TEXT (e.g. "246.13"), not numbers,
# even though they look numeric. Excel's COM layer auto-converts a bare
# numeric-looking string assigned via .Value into a real number, so those
# writes are done with a leading apostrophe (forces text entry, like typing
# it in the UI) and then the cell style is reset back to "Normal" so the
# quote-prefix formatting Excel applies doesn't leave a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice {
    param($Address, $NewValue)
    $ws.Range($Address).Value = "'" + $NewValue
    $ws.Range($Address).Style = "Normal"
}

# Column D ("Price") refreshes for otherwise-unchanged rows.
Set-TextPrice "D2"  "246.13"
Set-TextPrice "D4"  "5.367"
Set-TextPrice "D5"  "0.05868"
Set-TextPrice "D6"  "3.387"
Set-TextPrice "D7"  "6.376"
Set-TextPrice "D8"  "0.8133"
Set-TextPrice "D9"  "0.9815"
Set-TextPrice "D10" "0.1421"
Set-TextPrice "D11" "0.03640"
Set-TextPrice "D12" "0.07363"
Set-TextPrice "D13" "0.03013"
Set-TextPrice "D14" "4.469"
Set-TextPrice "D15" "0.09394"
Set-TextPrice "D16" "0.001597"
Set-TextPrice "D17" "0.04835"

# Row 18 (One / ONE): price refresh + it now carries the "Worst in 24h" tag.
Set-TextPrice "D18" "0.0005893"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextPrice "D19" "0.006226"
Set-TextPrice "D20" "0.004084"
Set-TextPrice "D21" "0.0009874"
Set-TextPrice "D22" "0.00009707"
Set-TextPrice "D23" "3.687"
Set-TextPrice "D26" "0.1296"
Set-TextPrice "D27" "0.0002472"
Set-TextPrice "D40" "0.03845"

# Rows 41-43 re-rank: KickToken moves up to 41, BKEXToken drops to 42,
# CEJI drops to 43 (and loses the "Worst in 24h" tag it used to carry as
# KickToken in row 43; the tag moved to row 18 above).
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice "D41" "0.006474"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.003002"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextPrice "D44" "0.005754"
Set-TextPrice "D45" "0.00005660"
Set-TextPrice "D46" "0.00000000750"
Set-TextPrice "D47" "0.6513"
Set-TextPrice "D48" "0.07787"
Set-TextPrice "D49" "0.00002101"
